$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.139.85'
$ws.Range('D3').Value = '2.362.56'
$ws.Range('E3').Value = '  +2.18%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''302.20'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.15%  '
$ws.Range('E6').Value = '  +0.87%  '
$ws.Range('E7').Value = '  -0.27%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  -2.80%  '
$ws.Range('E10').Value = '  -3.33%  '
$ws.Range('D11').Value = '''0.0801'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.26%  '
$ws.Range('D12').Value = '''7.16'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.30%  '
$ws.Range('E13').Value = '  -0.37%  '
$ws.Range('D14').Value = '2.723.04'
$ws.Range('E14').Value = '  +2.20%  '
$ws.Range('D15').Value = '2.355.50'
$ws.Range('E15').Value = '  +1.89%  '
$ws.Range('D16').Value = '''0.817'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.19%  '
$ws.Range('D17').Value = '''13.64'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.54%  '
$ws.Range('D18').Value = '46.106.09'
$ws.Range('E18').Value = '  -1.45%  '
$ws.Range('D19').Value = '''12.83'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.44%  '
$ws.Range('D20').Value = '0.0₃0968'
$ws.Range('E20').Value = '  +2.98%  '
$ws.Range('E21').Value = '  -0.93%  '
$ws.Range('D22').Value = '''67.53'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.13%  '
$ws.Range('D23').Value = '''246.46'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.67%  '
$ws.Range('E24').Value = '  -2.39%  '
$ws.Range('E25').Value = '  +0.12%  '
$ws.Range('E26').Value = '  -2.99%  '
$ws.Range('D27').Value = '''39.92'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.64%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').Value = '''9.82'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.39%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '''2.19'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.08%  '
$ws.Range('B30').Value = 'LidoDAOToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D30').Value = '''3.78'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +21.38%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '''21.04'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.31%  '
$ws.Range('D32').Value = '''2.79'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.27%  '
$ws.Range('E33').Value = '  -3.55%  '
$ws.Range('D34').Value = '''146.40'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.48%  '
$ws.Range('E35').Value = '  -2.60%  '
$ws.Range('E36').Value = '  -0.93%  '
$ws.Range('D37').Value = '''1.91'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.02%  '
$ws.Range('E38').Value = '  -1.81%  '
$ws.Range('D39').Value = '''15.09'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.51%  '
$ws.Range('D40').Value = '''3.97'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.98%  '
$ws.Range('D41').Value = '''0.0301'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.92%  '
$ws.Range('E42').Value = '  -4.75%  '
$ws.Range('D43').Value = '1.908.56'
$ws.Range('E43').Value = '  +3.57%  '
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').Value = '''92.27'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.25%  '
$ws.Range('E46').Value = '  -8.37%  '
$ws.Range('E47').Value = '  -6.38%  '
$ws.Range('D48').Value = '''8.27'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.79%  '
$ws.Range('D49').Value = '''97.97'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.81%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '''14.59'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.07%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.594.38'
$ws.Range('E51').Value = '  +2.08%  '
